$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-28 06:18:17'
$ws.Range('E3').Value = '2026-02-28 06:18:19'
$ws.Range('H3').Value = '85%'
$ws.Range('O3').Value = '-0.6 °C'
$ws.Range('E4').Value = '2026-02-28 06:18:22'
$ws.Range('M4').Value = '9.0 °C 5:51 TU'
$ws.Range('E5').Value = '2026-02-28 06:18:24'
$ws.Range('N5').Value = '-1.1 °C 5:51 TU'
$ws.Range('E6').Value = '2026-02-28 06:18:27'
$ws.Range('N6').Value = '9.9 °C 5:43 TU'
$ws.Range('E7').Value = '2026-02-28 06:18:29'
$ws.Range('M7').Value = '12.7 °C 5:34 TU'
$ws.Range('E8').Value = '2026-02-28 06:18:32'
$ws.Range('L8').Value = '26.3 km/h - 81º 5:49 TU'
$ws.Range('M8').Value = '8.9 °C 5:59 TU'
$ws.Range('E9').Value = '2026-02-28 06:18:35'
$ws.Range('N9').Value = '5.9 °C 5:55 TU'
$ws.Range('O9').Value = '7.3 °C'
$ws.Range('E10').Value = '2026-02-28 06:18:37'
$ws.Range('M10').Value = '9.4 °C 5:59 TU'
$ws.Range('O10').Value = '7.9 °C'
$ws.Range('E11').Value = '2026-02-28 06:18:39'
$ws.Range('N11').Value = '1.1 °C 5:58 TU'
$ws.Range('O11').Value = '3.2 °C'
$ws.Range('E12').Value = '2026-02-28 06:18:41'
$ws.Range('E13').Value = '2026-02-28 06:18:43'
$ws.Range('J13').Value = '1026.4 hPa'
$ws.Range('O13').Value = '1.0 °C'
$ws.Range('E14').Value = '2026-02-28 06:18:46'
$ws.Range('H14').Value = '97%'
$ws.Range('M14').Value = '12.5 °C 5:49 TU'
$ws.Range('O14').Value = '10.5 °C'
$ws.Range('E15').Value = '2026-02-28 06:18:48'
$ws.Range('E16').Value = '2026-02-28 06:18:50'
$ws.Range('H16').Value = '59%'
$ws.Range('L16').Value = '21.6 km/h - 216º 5:54 TU'
$ws.Range('N16').Value = '-2.1 °C 5:50 TU'
$ws.Range('E17').Value = '2026-02-28 06:18:52'
$ws.Range('H17').Value = '45%'
$ws.Range('K17').Value = '-0.1 MJ/m2'
$ws.Range('N17').Value = '2.5 °C 5:58 TU'
$ws.Range('O17').Value = '4.5 °C'
$ws.Range('E18').Value = '2026-02-28 06:18:55'
$ws.Range('M18').Value = '9.7 °C 5:50 TU'
$ws.Range('O18').Value = '8.4 °C'
$ws.Range('E19').Value = '2026-02-28 06:18:58'
$ws.Range('H19').Value = '67%'
$ws.Range('E20').Value = '2026-02-28 06:19:00'
$ws.Range('H20').Value = '38%'
$ws.Range('N20').Value = '-1.4 °C 5:54 TU'
$ws.Range('E21').Value = '2026-02-28 06:19:02'
$ws.Range('H21').Value = '77%'
$ws.Range('J21').Value = '1024.2 hPa'
$ws.Range('K21').Value = '-0.1 MJ/m2'
$ws.Range('N21').Value = '2.9 °C 5:54 TU'
$ws.Range('O21').Value = '4.9 °C'
$ws.Range('E22').Value = '2026-02-28 06:19:05'
$ws.Range('E23').Value = '2026-02-28 06:19:07'
$ws.Range('H23').Value = '66%'
$ws.Range('K23').Value = '-0.1 MJ/m2'
$ws.Range('E24').Value = '2026-02-28 06:19:10'
$ws.Range('J24').Value = '1023.6 hPa'
$ws.Range('O24').Value = '6.2 °C'
$ws.Range('E25').Value = '2026-02-28 06:19:12'
$ws.Range('H25').Value = '55%'
$ws.Range('O25').Value = '0.7 °C'
$ws.Range('E26').Value = '2026-02-28 06:19:14'
$ws.Range('J26').Value = '1024.0 hPa'
$ws.Range('N26').Value = '2.9 °C 5:57 TU'
$ws.Range('O26').Value = '4.4 °C'
$ws.Range('E27').Value = '2026-02-28 06:19:16'
$ws.Range('N27').Value = '0.4 °C 5:32 TU'
$ws.Range('O27').Value = '2.1 °C'
$ws.Range('E28').Value = '2026-02-28 06:19:19'
$ws.Range('E29').Value = '2026-02-28 06:19:22'
$ws.Range('H29').Value = '96%'
$ws.Range('O29').Value = '8.9 °C'
$ws.Range('E30').Value = '2026-02-28 06:19:25'
$ws.Range('M30').Value = '8.2 °C 5:58 TU'
$ws.Range('E31').Value = '2026-02-28 06:19:27'
$ws.Range('O31').Value = '10.0 °C'
$ws.Range('E32').Value = '2026-02-28 06:19:29'
$ws.Range('H32').Value = '91%'
$ws.Range('E33').Value = '2026-02-28 06:19:32'
$ws.Range('H33').Value = '71%'
$ws.Range('J33').Value = '1023.7 hPa'
$ws.Range('O33').Value = '4.9 °C'
$ws.Range('E34').Value = '2026-02-28 06:19:34'
$ws.Range('H34').Value = '65%'
$ws.Range('N34').Value = '-1.9 °C 5:57 TU'
$ws.Range('O34').Value = '-0.2 °C'
$ws.Range('E35').Value = '2026-02-28 06:19:37'
$ws.Range('J35').Value = '1023.0 hPa'
$ws.Range('N35').Value = '5.1 °C 5:59 TU'
$ws.Range('O35').Value = '6.5 °C'
$ws.Range('E36').Value = '2026-02-28 06:19:40'
$ws.Range('N36').Value = '8.2 °C 5:41 TU'
$ws.Range('O36').Value = '9.9 °C'
$ws.Range('E37').Value = '2026-02-28 06:19:42'
$ws.Range('N37').Value = '3.1 °C 5:43 TU'
$ws.Range('O37').Value = '4.5 °C'
$ws.Range('E38').Value = '2026-02-28 06:19:45'
$ws.Range('E39').Value = '2026-02-28 06:19:47'
$ws.Range('H39').Value = '47%'
$ws.Range('E40').Value = '2026-02-28 06:19:49'
$ws.Range('H40').Value = '94%'
$ws.Range('N40').Value = '1.8 °C 5:31 TU'
$ws.Range('O40').Value = '3.3 °C'
$ws.Range('E41').Value = '2026-02-28 06:19:51'
$ws.Range('J41').Value = '1023.2 hPa'
$ws.Range('O41').Value = '11.9 °C'
$ws.Range('E42').Value = '2026-02-28 06:19:54'
$ws.Range('E43').Value = '2026-02-28 06:19:57'
$ws.Range('H43').Value = '86%'
$ws.Range('N43').Value = '2.4 °C 5:59 TU'
$ws.Range('E44').Value = '2026-02-28 06:19:59'
$ws.Range('L44').Value = '20.9 km/h - 20º 5:52 TU'
$ws.Range('E45').Value = '2026-02-28 06:20:02'
$ws.Range('J45').Value = '1024.1 hPa'
$ws.Range('N45').Value = '4.9 °C 5:59 TU'
$ws.Range('O45').Value = '6.7 °C'
$ws.Range('E46').Value = '2026-02-28 06:20:04'
$ws.Range('J46').Value = '1023.2 hPa'
